$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 64
$prevRow = 63

# Column A holds a literal date-like text string (e.g. "2026/01/13"), not a
# real Excel date, matching every prior row in this log. Force a text
# number-format on the cell first so the value isn't auto-coerced into a
# date serial number when assigned.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/01/13"

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1144

# Match the cell formatting/style of the previous data row (63).
$ws.Range("A63:C63").Copy()
$ws.Range("A64:C64").PasteSpecial(-4122) # xlPasteFormats
